$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder: Rumania now ranks above Guatemala (rows 47/48 swap labels + data) ---
# Row 47 becomes Rumania with its freshly updated totals.
$ws.Cells.Item(47, 1).Value = "Rumania"
$ws.Cells.Item(47, 2).Value = 45902
$ws.Cells.Item(47, 3).Value = 1104
$ws.Cells.Item(47, 4).Value = 25794
$ws.Cells.Item(47, 5).Value = 17902
$ws.Cells.Item(47, 7).Value = 19
$ws.Cells.Item(47, 8).Value = 2206

# Row 48 becomes Guatemala, carrying the totals that used to belong to row 47.
$ws.Cells.Item(48, 1).Value = "Guatemala"
$ws.Cells.Item(48, 2).Value = 45053
$ws.Cells.Item(48, 4).Value = 31612
$ws.Cells.Item(48, 5).Value = 11707
$ws.Cells.Item(48, 8).Value = 1734

# --- Reorder: Groenlandia now ranks above Islas Malvinas (rows 210/211 swap labels) ---
# The underlying totals are identical between the two rows, only the labels move.
$ws.Cells.Item(210, 1).Value = "Groenlandia"
$ws.Cells.Item(211, 1).Value = "Islas Malvinas"

# --- Update "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Julio de 2020 a las 13:37"

# --- Updated country case counts ---
# Iran (row 14)
$ws.Cells.Item(14, 2).Value = 293606
$ws.Cells.Item(14, 3).Value = 2434
$ws.Cells.Item(14, 4).Value = 255144
$ws.Cells.Item(14, 5).Value = 22550
$ws.Cells.Item(14, 7).Value = 212
$ws.Cells.Item(14, 8).Value = 15912

# Kuwait (row 39)
$ws.Cells.Item(39, 2).Value = 64379
$ws.Cells.Item(39, 3).Value = 606
$ws.Cells.Item(39, 4).Value = 55057
$ws.Cells.Item(39, 5).Value = 8884
$ws.Cells.Item(39, 7).Value = 5
$ws.Cells.Item(39, 8).Value = 438

# Barein (row 52)
$ws.Cells.Item(52, 5).Value = 3301
$ws.Cells.Item(52, 7).Value = 1
$ws.Cells.Item(52, 8).Value = 141

# Suiza (row 55)
$ws.Cells.Item(55, 2).Value = 34477
$ws.Cells.Item(55, 3).Value = 65
$ws.Cells.Item(55, 5).Value = 1799
$ws.Cells.Item(55, 7).Value = 1
$ws.Cells.Item(55, 8).Value = 1978

# Nepal (row 67)
$ws.Cells.Item(67, 2).Value = 18752
$ws.Cells.Item(67, 3).Value = 139
$ws.Cells.Item(67, 4).Value = 13754
$ws.Cells.Item(67, 5).Value = 4950
$ws.Cells.Item(67, 7).Value = 3
$ws.Cells.Item(67, 8).Value = 48

# Madagascar (row 85)
$ws.Cells.Item(85, 2).Value = 9690
$ws.Cells.Item(85, 3).Value = 395
$ws.Cells.Item(85, 4).Value = 6260
$ws.Cells.Item(85, 5).Value = 3339
$ws.Cells.Item(85, 7).Value = 6
$ws.Cells.Item(85, 8).Value = 91

# Burkina Faso (row 146)
$ws.Cells.Item(146, 2).Value = 1100
$ws.Cells.Item(146, 3).Value = 14
$ws.Cells.Item(146, 4).Value = 926
$ws.Cells.Item(146, 5).Value = 121

# Vietnam (row 163)
$ws.Cells.Item(163, 2).Value = 423
$ws.Cells.Item(163, 3).Value = 3
$ws.Cells.Item(163, 5).Value = 58
